# Updated serpent xsdir paths
# Append ".serp" to the Serpent xsdir path entries (column E) on the
# "Libraries" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Libraries")

$ws.Range("E2").Value = "/home/mcnp/xs/xsdir_mcnp6.2.serp"
$ws.Range("E3").Value = "/home/mcnp/xs/xsdir_mcnp6.2_old.serp"
$ws.Range("E4").Value = "/home/mcnp/xs/xsdir_mcnp6.2_jeff33_endfb71_fendl32b_irdff105_tt.serp"
$ws.Range("E5").Value = "/home/mcnp/xs/xsdir_mcnp6.2_fendl32b_rw.serp"
$ws.Range("E6").Value = "/home/mcnp/xs/xsdir_mcnp6.2.serp"
$ws.Range("E7").Value = "/home/mcnp/xs/xsdir_mcnp6.2_endfb8.serp"

$ws.Range("E7").Select()
